$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell B45 (week 44) from 63 to 444
$ws.Range("B45").Value = 444

# Add new row 46 for week 45 of 2024
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = 580
